# "Generate Report for Handback" -- the handback report generator reran
# and refreshed the handoff/handback timestamps for the file that was just
# regenerated (63bdf7ad-...), while the already-handed-back file
# (f4ef5a11-...) keeps its prior timestamps untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" for 63bdf7ad-....md moves forward.
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("G2").Value = "2016-08-13 19:00:26"

# --- zh-cn sheet ---
# Correspond Handoff Datetime / Correspond Handback DateTime for the
# 63bdf7ad-....md row move forward.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H2").Value = "2016-08-13 19:00:18"
$wsZh.Range("K2").Value = "2016-08-13 19:00:47"

# --- de-de sheet ---
# Correspond Handoff Datetime / Correspond Handback DateTime for the
# 63bdf7ad-....md row move forward.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H2").Value = "2016-08-13 19:00:26"
$wsDe.Range("K2").Value = "2016-08-13 19:00:57"
